$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.781.47'
$ws.Range("E2").Value = '  -2.60%  '

# Row 3
$ws.Range("D3").Value = '1.559.75'
$ws.Range("E3").Value = '  -0.57%  '

# Row 4
$ws.Range("E4").Value = '  +0.24%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '205.94'
$ws.Range("E5").Value = '  -1.06%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.486'
$ws.Range("E6").Value = '  -2.23%  '

# Row 7
$ws.Range("E7").Value = '  +0.29%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.90'
$ws.Range("E8").Value = '  -0.51%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.247'
$ws.Range("E9").Value = '  -0.50%  '

# Row 10
$ws.Range("E10").Value = '  -1.35%  '

# Row 11
$ws.Range("E11").Value = '  -0.15%  '

# Row 12
$ws.Range("D12").Value = '1.782.99'
$ws.Range("E12").Value = '  -0.49%  '

# Row 13
$ws.Range("D13").Value = '1.571.01'
$ws.Range("E13").Value = '  +0.10%  '

# Row 14
$ws.Range("E14").Value = '  -2.33%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.511'
$ws.Range("E15").Value = '  -1.59%  '

# Row 16
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '26.818.85'
$ws.Range("E16").Value = '  -2.36%  '

# Row 17
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.54'
$ws.Range("E17").Value = '  -2.87%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '214.26'
$ws.Range("E18").Value = '  +0.03%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.30'
$ws.Range("E19").Value = '  +0.34%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0677'
$ws.Range("E20").Value = '  -1.79%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.08'
$ws.Range("E22").Value = '  -1.04%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.32'
$ws.Range("E23").Value = '  -2.49%  '

# Row 24
$ws.Range("E24").Value = '  -1.27%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.66'
$ws.Range("E25").Value = '  -0.99%  '

# Row 26
$ws.Range("E26").Value = '  -0.87%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.83'
$ws.Range("E27").Value = '  -1.29%  '

# Row 28
$ws.Range("E28").Value = '  +0.25%  '

# Row 29
$ws.Range("E29").Value = '  -1.65%  '

# Row 30
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.11'
$ws.Range("E30").Value = '  -3.95%  '

# Row 31
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0461'
$ws.Range("E31").Value = '  -2.06%  '

# Row 32
$ws.Range("E32").Value = '  -1.79%  '

# Row 33
$ws.Range("D33").Value = '1.387.11'
$ws.Range("E33").Value = '  +1.83%  '

# Row 34
$ws.Range("E34").Value = '  -1.49%  '

# Row 35
$ws.Range("E35").Value = '  +0.62%  '

# Row 36
$ws.Range("E36").Value = '  -0.21%  '

# Row 37
$ws.Range("E37").Value = '  -4.61%  '

# Row 38
$ws.Range("E38").Value = '  -3.04%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.806'
$ws.Range("E39").Value = '  -1.93%  '

# Row 40
$ws.Range("E40").Value = '  -3.98%  '

# Row 41
$ws.Range("E41").Value = '  +0.25%  '

# Row 42
$ws.Range("E42").Value = '  +2.91%  '

# Row 43
$ws.Range("E43").Value = '  +2.64%  '

# Row 44
$ws.Range("E44").Value = '  +1.65%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.76'
$ws.Range("E45").Value = '  -1.83%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '63.05'
$ws.Range("E46").Value = '  -1.80%  '

# Row 47
$ws.Range("D47").Value = '1.695.85'
$ws.Range("E47").Value = '  -0.56%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.56'
$ws.Range("E48").Value = '  +0.02%  '

# Row 49
$ws.Range("D49").Value = '0.0₇0984'
$ws.Range("E49").Value = '  -0.37%  '

# Row 50
$ws.Range("E50").Value = '  -0.39%  '

# Row 51
$ws.Range("E51").Value = '  -1.19%  '
